# Updates the cryptos list Price (D) and Volume(1h) (E) columns
# to the latest scraped values. A leading apostrophe forces Excel
# to store the value as literal text (preserving things like
# trailing zeros / thousands-dot formatting); Style is reset back
# to "Normal" afterwards so no extra number-format style sticks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.398.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -8.05%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.676.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -7.00%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.41%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'215.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -6.67%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.007"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.30%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4939"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -17.07%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.2589"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -6.76%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'21.61"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -7.49%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.06116"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -10.46%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07282"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -3.69%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.717.83"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -4.72%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'4.401"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -7.76%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.5701"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -8.69%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'1.904.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -7.01%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.000008136"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -12.79%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'64.00"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -15.17%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'26.429.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -7.86%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'4.965"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -9.37%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +0.32%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'10.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -6.90%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'182.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -13.57%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'6.146"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -10.36%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +0.34%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'144.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -6.75%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'7.470"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -4.82%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.1124"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -11.93%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'15.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -5.22%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.308"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -8.60%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.05647"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -8.22%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.316"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -7.34%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.459"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -8.68%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'3.442"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -8.21%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.618"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -6.06%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.9980"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -6.08%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D37").Value = "'0.5846"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -8.62%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'2.626"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -3.17%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.01574"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -8.05%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'1.065.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -6.00%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'5.866"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -9.14%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.8462"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.78%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'1.001"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.39%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'97.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.12%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.835.73"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -6.40%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'55.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -7.63%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -6.16%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.005"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.17%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'8.063"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -3.30%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.4326"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.69%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.05179"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -5.34%  "
$ws.Range("E51").Style = "Normal"
